$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "م" (row number) becomes 1
$ws.Range("A7").Value = 1

# Row 7: product name -> text (style numFmt changes General -> Text for this merged block)
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "سائل ريد"

# Row 7: stock/balance -> text (style numFmt changes General -> Text for this merged block)
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "12:0"

# Row 7: reorder limit -> stored as text but keeps its existing numeric display format
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "0"
$ws.Range("L7").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

# Row 7: price -> text (reuses the same style as product name block)
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "100.00"

# Row 7: sell price -> stored as text but keeps its existing numeric display format
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "100.0000"
$ws.Range("P7").NumberFormat = "0.00"

# Row 7: transaction count -> text
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"

# Row 8: numeric total
$ws.Range("P8").Value = 100

# Row 9: refreshed generation timestamp
$ws.Range("A9").Value = "Saturday, 26 July, 2025 12:08 AM"
